# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" sheet and the "全部类型" sheet, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 468
$ws1.Range("F3").Value = 5616
$ws1.Range("F5").Value = 70
$ws1.Range("F10").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 468
$ws4.Range("F3").Value = 5616
$ws4.Range("F6").Value = 70
$ws4.Range("F12").Value = 23
